$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Risk_Tracking_Log")
$ws.Activate()

# --- Row 15: new risk entry (Application Metabase and vue.js do not meet project needs) ---
$ws.Range("A15").Value = "6"
$ws.Range("B15").Value = "Open"
$ws.Range("C15").Value = "Medium"
$ws.Range("D15").Value = "Medium"
$ws.Range("F15").Value = "Application Metabase and vue.js do not meet project needs"
$ws.Range("G15").Value = "New applications will need to be sourced"
$ws.Range("H15").Value = "Schedule" + [char]10 + "Technology" + [char]10 + "Reliability of systems"
$ws.Range("I15").Value = "Unable to execute required functions, such as visualisations"
$ws.Range("J15").Value = "Difficulty in integrating software into project"
$ws.Range("K15").Value = "Mitigation"
$ws.Range("L15").Value = "Team members report issues immediately the arise"
$ws.Range("M15").Value = "When indication start to appear that software is not compatible or doesn't work as desired, investigation starts urgently for new software to take its place"

# --- Row heights: rows 13 and 15 grow to fit wrapped text (45 -> 56.25) ---
$ws.Rows.Item(13).RowHeight = 56.25
$ws.Rows.Item(15).RowHeight = 56.25

# --- Selection / scroll position bookkeeping ---
$ws.Range("A13").Select()
